$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-hide the previously filtered-out rows and clear the column filter criteria
[void]$ws.ShowAllData()

# Append the new feature row
$ws.Range("A13").Value2 = "Reroll below threshold"
$ws.Range("B13").Value2 = "Re-roll dice that are below a certain threshold, like 1 & 2 for great weapon master"
$ws.Range("D13").Value2 = "Sean Steele - ssteele1812@gmail.com"

# Move the active selection to match the author's final cursor position
[void]$ws.Range("C13").Select()
